# [#89] Add addressRegion (County) to rollUp
#
# Insert a new "Recipient Org:County" column into the "grants" summary sheet,
# directly after the existing "Recipient Org:City" column, pushing every
# column from "Recipient Org:Country" onward one slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants")

# "Recipient Org:City" lives in column R (18); the new column goes in S (19),
# shifting the former S:AL ("Recipient Org:Country" .. "Data Source") to T:AM.
$ws.Columns.Item(19).Insert()
$ws.Cells.Item(1, 19).Value = "Recipient Org:County"
